# Applies the czech-republic cfl-group-a 2023-2024 update:
#  - Swap the (previously mis-ordered) row pairs so that match data lines
#    up with the correct date-of-match ordering (columns F:V only; the
#    leading Indice/pais/torneio/temporada/data_partida columns A:E are
#    identical within each swapped pair so they are left untouched).
#  - Append three brand-new fixtures as rows 58-60.
#  - Extend the sheet's used-range dimension to A1:V60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowDataFV {
    param($sheet, [int]$rowA, [int]$rowB)

    # Columns F (6) through V (22). Per-cell Value2 read (plain ".Value"
    # reads resolve to the property descriptor rather than the cell's
    # contents in this host) + per-cell Value write.
    for ($col = 6; $col -le 22; $col++) {
        $cellA = $sheet.Cells.Item($rowA, $col)
        $cellB = $sheet.Cells.Item($rowB, $col)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# Row pairs whose F:V content needs to be swapped
Swap-RowDataFV $ws 10 11
Swap-RowDataFV $ws 14 15
Swap-RowDataFV $ws 16 17
Swap-RowDataFV $ws 30 31
Swap-RowDataFV $ws 48 49
Swap-RowDataFV $ws 54 55

# Append three new fixture rows (58, 59, 60), cloning row 57's formatting
# first so number formats / styles (Indice bold+border, date format, etc.)
# stay consistent with the rest of the sheet.
$ws.Range("A57:V57").Copy($ws.Range("A58:V58"))
$ws.Range("A57:V57").Copy($ws.Range("A59:V59"))
$ws.Range("A57:V57").Copy($ws.Range("A60:V60"))

# Row 58 -> Indice 57: Loko Vltavin vs Povltavska FA
$row = 58
$ws.Cells.Item($row, 1).Value = 57
$ws.Cells.Item($row, 2).Value = "czech-republic"
$ws.Cells.Item($row, 3).Value = "cfl-group-a"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45192.42708333334
$ws.Cells.Item($row, 6).Value = "Vltavin"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Povltavska FA"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.75
$ws.Cells.Item($row, 11).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 12).Value = 1.56
$ws.Cells.Item($row, 13).Value = "23/09/2023 10:06"
$ws.Cells.Item($row, 14).Value = 3.74
$ws.Cells.Item($row, 15).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 16).Value = 4.27
$ws.Cells.Item($row, 17).Value = "23/09/2023 10:06"
$ws.Cells.Item($row, 18).Value = 3.44
$ws.Cells.Item($row, 19).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 20).Value = 4.74
$ws.Cells.Item($row, 21).Value = "23/09/2023 10:06"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/loko-vltavin-povltavska-fa/A5MLXJJj/"

# Row 59 -> Indice 58: Ceske Budejovice B vs Pisek
$row = 59
$ws.Cells.Item($row, 1).Value = 58
$ws.Cells.Item($row, 2).Value = "czech-republic"
$ws.Cells.Item($row, 3).Value = "cfl-group-a"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45192.4375
$ws.Cells.Item($row, 6).Value = "Ceske Budejovice B"
$ws.Cells.Item($row, 7).Value = 4
$ws.Cells.Item($row, 8).Value = "Pisek"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.83
$ws.Cells.Item($row, 11).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 12).Value = 1.59
$ws.Cells.Item($row, 13).Value = "23/09/2023 10:21"
$ws.Cells.Item($row, 14).Value = 3.75
$ws.Cells.Item($row, 15).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 16).Value = 4.4
$ws.Cells.Item($row, 17).Value = "23/09/2023 10:21"
$ws.Cells.Item($row, 18).Value = 3.15
$ws.Cells.Item($row, 19).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 20).Value = 4.35
$ws.Cells.Item($row, 21).Value = "23/09/2023 10:21"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/ceske-budejovice-pisek/QX6KhLeG/"

# Row 60 -> Indice 59: Kraluv Dvur vs FK Robstav
$row = 60
$ws.Cells.Item($row, 1).Value = 59
$ws.Cells.Item($row, 2).Value = "czech-republic"
$ws.Cells.Item($row, 3).Value = "cfl-group-a"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45192.4375
$ws.Cells.Item($row, 6).Value = "Kraluv Dvur"
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = "FK Robstav"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 3.36
$ws.Cells.Item($row, 11).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 12).Value = 3.19
$ws.Cells.Item($row, 13).Value = "23/09/2023 10:26"
$ws.Cells.Item($row, 14).Value = 3.79
$ws.Cells.Item($row, 15).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 16).Value = 4.02
$ws.Cells.Item($row, 17).Value = "23/09/2023 10:26"
$ws.Cells.Item($row, 18).Value = 1.76
$ws.Cells.Item($row, 19).Value = "21/09/2023 21:42"
$ws.Cells.Item($row, 20).Value = 1.91
$ws.Cells.Item($row, 21).Value = "23/09/2023 10:26"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/kraluv-dvur-fk-robstav/I73CfsR3/"
